$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date-range header text (merged A4:I4)
$ws.Range("A4").Value = "Từ ngày 22/07/2018 tới ngày 22/08/2018"

# Delete row 7 (the "11/07/2018 / dsadwa / dwasdwada / Tin học..." entry) -
# remaining rows 8-12 shift up to become rows 7-11
$ws.Rows.Item(7).Delete()

# Renumber the STT column (A) for the now-shifted data rows 7-11 -> 1..5
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
